# repull data, push all data, mean calculation
# Updates the "dSF" column (F) values for the affected rows to reflect
# the re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 7
    3  = 1
    4  = -1
    6  = -2
    7  = 4
    9  = -5
    14 = -1
    16 = 5
    17 = -4
    18 = -5
    19 = -4
    23 = -2
    24 = -3
    26 = -1
    31 = -7
    32 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
